$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HSI review")

# Fix typo in C15: "requireemnts" -> "requirements"
$ws.Range("C15").Value = "In Software context: " + [char]10 + "- The inputs and outputs in the requirements shall be the same as output signals in sotware context for ex: " + [char]10 + "in SRS_018 I shall see ""Tail LEDs status"" in the context instead it's not exist and all the rest signals like that " + [char]10 + "- For the input side the signal for ex: ""Mode signal"" is input to ""input feature block"" and output from it at the same time !"

# Set Acceptance decision for points 003-010 (row13) and Software context point (row15) to Accepted
$ws.Range("D13").Value = "Accepted"
$ws.Range("D15").Value = "Accepted"

# Update active selection to D15
$ws.Range("D15").Select()
